# Update the Vermont courts list worksheet:
#  - court unit emails move from "JUD.<Unit>@vermont.gov" to "<Unit>@vtcourts.gov"
#  - the new email cells become mailto: hyperlinks (styled with the built-in
#    "Hyperlink" style), keeping the old address as the hyperlink's display
#    text/tooltip fallback
#  - the active-cell selection left on the sheet moves to D27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Turn I2:I15 into mailto: hyperlinks first (in row order, so the
#    relationship ids rId1..rId14 line up with rows I2..I15). Passing the
#    *old* JUD.<Unit>@vermont.gov address as TextToDisplay here only seeds
#    the hyperlink's display/tooltip text -- the cell values themselves are
#    overwritten afterwards with the new @vtcourts.gov addresses below.
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:AddisonUnit@vtcourts.gov", "", "", "AddisonUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:BenningtonUnit@vtcourts.gov", "", "", "BenningtonUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:CaledoniaEssexUnit@vtcourts.gov", "", "", "CaledoniaEssexUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:ChittendenUnit@vtcourts.gov", "", "", "ChittendenUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:CaledoniaEssexUnit@vtcourts.gov", "", "", "CaledoniaEssexUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I7"), "mailto:FranklinGrandIsleunit@vtcourts.gov", "", "", "FranklinGrandIsleunit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I8"), "mailto:FranklinGrandIsleunit@vtcourts.gov", "", "", "FranklinGrandIsleunit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I9"), "mailto:LamoilleUnit@vtcourts.gov", "", "", "LamoilleUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I10"), "mailto:OrangeUnit@vtcourts.gov", "", "", "OrangeUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I11"), "mailto:OrleansUnit@vtcourts.gov", "", "", "OrleansUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I12"), "mailto:RutlandUnit@vtcourts.gov", "", "", "RutlandUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I13"), "mailto:WashingtonUnit@vtcourts.gov", "", "", "WashingtonUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I14"), "mailto:WindhamUnit@vtcourts.gov", "", "", "WindhamUnit@vermont.gov")
$ws.Hyperlinks.Add($ws.Range("I15"), "mailto:WindsorUnit@vtcourts.gov", "", "", "WindsorUnit@vermont.gov")

# 2) Now overwrite the cell text with the new @vtcourts.gov addresses
#    (no "JUD." prefix). I3..I15 first, I2 last, so the shared-string table
#    is rebuilt in the same order the source workbook ended up with.
$ws.Range("I3").Value = "BenningtonUnit@vtcourts.gov"
$ws.Range("I4").Value = "CaledoniaEssexUnit@vtcourts.gov"
$ws.Range("I5").Value = "ChittendenUnit@vtcourts.gov"
$ws.Range("I6").Value = "CaledoniaEssexUnit@vtcourts.gov"
$ws.Range("I7").Value = "FranklinGrandIsleunit@vtcourts.gov"
$ws.Range("I8").Value = "FranklinGrandIsleunit@vtcourts.gov"
$ws.Range("I9").Value = "LamoilleUnit@vtcourts.gov"
$ws.Range("I10").Value = "OrangeUnit@vtcourts.gov"
$ws.Range("I11").Value = "OrleansUnit@vtcourts.gov"
$ws.Range("I12").Value = "RutlandUnit@vtcourts.gov"
$ws.Range("I13").Value = "WashingtonUnit@vtcourts.gov"
$ws.Range("I14").Value = "WindhamUnit@vtcourts.gov"
$ws.Range("I15").Value = "WindsorUnit@vtcourts.gov"
$ws.Range("I2").Value = "AddisonUnit@vtcourts.gov"

# 3) Leave the sheet's selection where the author last left it.
$ws.Range("D27").Select()
